# Localization strings: add FORMATTED_UI_GAME_CURRENT entry at the end of
# the table (row 52) so the "current total" can be shown on the bridge.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- add the new localization row --------------------------------------
$ws.Range("A52").Value = "FORMATTED_UI_GAME_CURRENT"
$ws.Range("B52").Value = "Current Total: {0}"
$ws.Range("C52").Value = "XXXX"
$ws.Range("D52").Value = "XXXX"
$ws.Range("E52").Value = "XXXX"

# match the formatting used by the rest of the table (style ids 3 / 4:
# left/top aligned, column A without wrap, columns B:E with wrap)
$ws.Range("A52").HorizontalAlignment = -4131
$ws.Range("A52").VerticalAlignment = -4160
$ws.Range("A52").WrapText = $false
$ws.Range("B52:E52").HorizontalAlignment = -4131
$ws.Range("B52:E52").VerticalAlignment = -4160
$ws.Range("B52:E52").WrapText = $true

# --- duplicate the conditional-format rule's dxf, the way Excel does
#     when a rule is re-saved, and repoint the surviving rule at it -----
$rng = $ws.Range("A1:XFD1048576")
$fc1 = $rng.FormatConditions.Item(1)
$origColor = $fc1.Interior.Color
$fc2 = $rng.FormatConditions.Add(1, 3, $fc1.Formula1)
$fc2.Interior.Color = $origColor
$fc1.Delete()

# --- scroll / selection bookkeeping so the view matches what a user
#     would see after scrolling down to add the new row ----------------
$ws.Activate()
$excel.Goto($ws.Range("A19"), $true)
$ws.Range("E52").Select()
